# Performer A3 - Sentiment Analysis : Done
# Rename the worksheet to reflect its contents
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "TargetChannelDetails"

# Add a new target channel row (Fireship) - enter the channel ID before
# the channel name so the shared-string table order matches
$ws.Range("B3").Value = "UCsBjURrPoezykLs9EqgamOA"
$ws.Range("A3").Value = "Fireship"
$ws.Range("C3").Value = "yt_queue"

# Match B3's font to B2's (bold, black) font
$ws.Range("B3").Font.Color = $ws.Range("B2").Font.Color

# Nudge the font on the data cells (A2, C2, A3, C3) to pick up the
# same "applied" font formatting used across the sheet
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("C2").Font.ThemeColor = 1
$ws.Range("A3").Font.ThemeColor = 1
$ws.Range("C3").Font.ThemeColor = 1

# Move the active selection to A4, matching the saved view state
[void]$ws.Range("A4").Select()
